$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "Distribution channel code" column is being inserted right before the
# existing "Actual" column (I), pushing "Actual" one column over to J.
$ws.Columns("I:I").Insert()

# Give the newly inserted column its own width (matches the authored layout).
$ws.Columns("I:I").ColumnWidth = 21.666666666666668

# Header text for the new column.
$ws.Range("I1").Value = "Distribution channel code"

# Data rows for the new column.
$ws.Range("I2").Value = "TR"
$ws.Range("I3").Value = "GO"
